# Generate Report for Handback
#
# - Status moves from "Ready for handoff" to "Handed back: in sync with en-US"
#   (text used throughout Overview / zh-cn / de-de sheets).
# - Each language sheet gains a "Latest Target File" (F) and "Latest Handback
#   File" (G) hyperlink, mirroring the source-file / target-file links already
#   present in columns A and D.
# - "Latest Handback DateTime" (H) is stamped with the handback time for each
#   language.

$wb = $excel.ActiveWorkbook

$oldStatus = "Ready for handoff"
$newStatus = "Handed back: in sync with en-US"

function Get-HyperlinkAddress($ws, $addr) {
    foreach ($hl in $ws.Hyperlinks) {
        if ($hl.Range.Address() -eq $addr) {
            return $hl.Address()
        }
    }
    return $null
}

function Set-LikeHyperlink($ws, $targetCellAddr, $likeCellAddr, $displayText) {
    # Re-creates, on $targetCellAddr, a hyperlink pointing at the same
    # address as the existing hyperlink found on $likeCellAddr, using the
    # same visible text, and applies the same (blue/underlined) look used
    # by the workbook's existing hyperlink cells.
    $url = Get-HyperlinkAddress $ws $likeCellAddr
    $cell = $ws.Range($targetCellAddr)
    $ws.Hyperlinks.Add($cell, $url, [Type]::Missing, [Type]::Missing, $displayText) | Out-Null
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 11
    $cell.Font.Underline = 2
    $cell.Font.Color = 15570276
}

# ---------------------------------------------------------------------
# 1. Status column: "Ready for handoff" -> "Handed back: in sync with en-US"
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
foreach ($addr in @("B2", "C2", "B3", "C3")) {
    $r = $overview.Range($addr)
    if ($r.Value2() -eq $oldStatus) {
        $r.Value2 = $newStatus
    }
}

foreach ($sheetName in @("zh-cn", "de-de")) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in @("C2", "C3")) {
        $r = $ws.Range($addr)
        if ($r.Value2() -eq $oldStatus) {
            $r.Value2 = $newStatus
        }
    }
}

# ---------------------------------------------------------------------
# 2. zh-cn sheet: populate F/G (target + handback file) and stamp H
#    (handback datetime) for both data rows.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

Set-LikeHyperlink $zhcn "F2" "`$A`$2" $zhcn.Range("A2").Value2()
Set-LikeHyperlink $zhcn "G2" "`$D`$2" $zhcn.Range("D2").Value2()
Set-LikeHyperlink $zhcn "F3" "`$A`$3" $zhcn.Range("A3").Value2()
Set-LikeHyperlink $zhcn "G3" "`$D`$3" $zhcn.Range("D3").Value2()

$zhcn.Range("H2").Value2 = "2016-03-12 06:32:31"
$zhcn.Range("H3").Value2 = "2016-03-12 06:32:31"

# ---------------------------------------------------------------------
# 3. de-de sheet: same treatment, but with its own handback timestamp.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

Set-LikeHyperlink $dede "F2" "`$A`$2" $dede.Range("A2").Value2()
Set-LikeHyperlink $dede "G2" "`$D`$2" $dede.Range("D2").Value2()
Set-LikeHyperlink $dede "F3" "`$A`$3" $dede.Range("A3").Value2()
Set-LikeHyperlink $dede "G3" "`$D`$3" $dede.Range("D3").Value2()

$dede.Range("H2").Value2 = "2016-03-12 06:32:37"
$dede.Range("H3").Value2 = "2016-03-12 06:32:37"
